$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.339.57'
$ws.Range("E2").Value = '  +8.65%  '

$ws.Range("D3").Value = '1.597.05'
$ws.Range("E3").Value = '  +7.92%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.28%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9948'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '288.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3691'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.75%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3395'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.75'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.137'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.20%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07033'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.63%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.65'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.917'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.622'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.87%  '

$ws.Range("B16").Value = 'Dai'
$ws.Range("C16").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9952'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.76%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.596.36'
$ws.Range("E17").Value = '  +7.77%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001078'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06598'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +11.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +12.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +11.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.021'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.96%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.36%  '

$ws.Range("D24").Value = '22.372.45'
$ws.Range("E24").Value = '  +8.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.394'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.34%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.500'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +15.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.41%  '

$ws.Range("D29").Value = '1.774.82'
$ws.Range("E29").Value = '  +8.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.54'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.96%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.175'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.051'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +21.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9458'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +15.92%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08229'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.40%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.591'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.285'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +11.79%  '

$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +12.31%  '

$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.596'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06114'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.237'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.27%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02208'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.87%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2025'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9944'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.76%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5886'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.52%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.670'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.99%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5686'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.963'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06816'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.62%  '
